$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 40093750
$ws.Range("I132").Value = 43738452
$ws.Range("J132").Value = 2033.3334
$ws.Range("K132").Value = 131215356
$ws.Range("L132").Value = 6100.0002
$ws.Range("M132").Value = -131212826
$ws.Range("N132").Value = -11160.0002
$ws.Range("H137").Value = 478933.34
$ws.Range("J137").Value = 84567.164
$ws.Range("L137").Value = 253701.492
$ws.Range("N137").Value = -258801.492

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 670466
$ws.Range("I32").Value = 3488.9744
$ws.Range("K32").Value = 3488.9744
$ws.Range("M32").Value = -3201.9744
$ws.Range("H61").Value = 8761.809999999999
$ws.Range("I61").Value = 11832.154
$ws.Range("J61").Value = 3772.5
$ws.Range("K61").Value = 11832.154
$ws.Range("L61").Value = 3772.5
$ws.Range("M61").Value = -11620.154
$ws.Range("N61").Value = -4196.5
$ws.Range("H74").Value = 3834.1428
$ws.Range("I74").Value = 621.0909
$ws.Range("K74").Value = 621.0909
$ws.Range("M74").Value = 252.9091
$ws.Range("H77").Value = 3834.1428
$ws.Range("I77").Value = 621.0909
$ws.Range("K77").Value = 3105.4545
$ws.Range("M77").Value = 1262.5455
$ws.Range("H132").Value = 2779288.2
$ws.Range("I132").Value = 4465153.5
$ws.Range("J132").Value = 2568.7646
$ws.Range("K132").Value = 13395460.5
$ws.Range("L132").Value = 7706.293799999999
$ws.Range("M132").Value = -13392930.5
$ws.Range("N132").Value = -12766.2938
$ws.Range("H136").Value = 8761.809999999999
$ws.Range("I136").Value = 11832.154
$ws.Range("J136").Value = 3772.5
$ws.Range("K136").Value = 35496.462
$ws.Range("L136").Value = 11317.5
$ws.Range("M136").Value = -32946.462
$ws.Range("N136").Value = -16417.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 4705
$ws.Range("I26").Value = 4705
$ws.Range("K26").Value = 4705
$ws.Range("M26").Value = -4413
$ws.Range("H75").Value = 6169
$ws.Range("I75").Value = 2338
$ws.Range("J75").Value = 10000
$ws.Range("K75").Value = 2338
$ws.Range("L75").Value = 10000
$ws.Range("M75").Value = -1402
$ws.Range("N75").Value = -11872
$ws.Range("H78").Value = 6169
$ws.Range("I78").Value = 2338
$ws.Range("J78").Value = 10000
$ws.Range("K78").Value = 7014
$ws.Range("L78").Value = 30000
$ws.Range("M78").Value = -2334
$ws.Range("N78").Value = -39360
$ws.Range("H96").Value = 4276
$ws.Range("I96").Value = 1414
$ws.Range("J96").Value = 10000
$ws.Range("K96").Value = 1414
$ws.Range("L96").Value = 10000
$ws.Range("M96").Value = 1332
$ws.Range("N96").Value = -15492
$ws.Range("H134").Value = 18548554
$ws.Range("I134").Value = 30350642
$ws.Range("J134").Value = 2416.2856
$ws.Range("K134").Value = 91051926
$ws.Range("L134").Value = 7248.8568
$ws.Range("M134").Value = -91049391
$ws.Range("N134").Value = -12318.8568

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 16351.85
$ws.Range("I31").Value = 9293.666999999999
$ws.Range("J31").Value = 26939.125
$ws.Range("K31").Value = 9293.666999999999
$ws.Range("L31").Value = 26939.125
$ws.Range("M31").Value = -8998.666999999999
$ws.Range("N31").Value = -27529.125
$ws.Range("H34").Value = 16351.85
$ws.Range("I34").Value = 9293.666999999999
$ws.Range("J34").Value = 26939.125
$ws.Range("K34").Value = 9293.666999999999
$ws.Range("L34").Value = 26939.125
$ws.Range("M34").Value = -9091.666999999999
$ws.Range("N34").Value = -27343.125
$ws.Range("H58").Value = 5110406.5
$ws.Range("I58").Value = 7937622
$ws.Range("J58").Value = 21418.8
$ws.Range("K58").Value = 7937622
$ws.Range("L58").Value = 21418.8
$ws.Range("M58").Value = -7937419
$ws.Range("N58").Value = -21824.8
$ws.Range("H103").Value = 11158.667
$ws.Range("I103").Value = 6212
$ws.Range("J103").Value = 21052
$ws.Range("K103").Value = 6212
$ws.Range("L103").Value = 21052
$ws.Range("M103").Value = -5040
$ws.Range("N103").Value = -23396
$ws.Range("H132").Value = 22231542
$ws.Range("I132").Value = 166667680
$ws.Range("K132").Value = 500003040
$ws.Range("M132").Value = -500000510
$ws.Range("H134").Value = 10778053
$ws.Range("I134").Value = 13160101
$ws.Range("J134").Value = 6252161.5
$ws.Range("K134").Value = 39480303
$ws.Range("L134").Value = 18756484.5
$ws.Range("M134").Value = -39477768
$ws.Range("N134").Value = -18761554.5
$ws.Range("H136").Value = 5110406.5
$ws.Range("I136").Value = 7937622
$ws.Range("J136").Value = 21418.8
$ws.Range("K136").Value = 23812866
$ws.Range("L136").Value = 64256.39999999999
$ws.Range("M136").Value = -23810316
$ws.Range("N136").Value = -69356.39999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 104.07143
$ws.Range("I12").Value = 185
$ws.Range("J12").Value = 82
$ws.Range("K12").Value = 555
$ws.Range("L12").Value = 246
$ws.Range("M12").Value = -382
$ws.Range("N12").Value = -592

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = $null
$ws.Range("N96").Value = 0
$ws.Range("H107").Value = 62735.562
$ws.Range("I107").Value = 66877.92999999999
$ws.Range("J107").Value = 600
$ws.Range("K107").Value = 66877.92999999999
$ws.Range("L107").Value = 600
$ws.Range("M107").Value = -64957.92999999999
$ws.Range("N107").Value = -4440
$ws.Range("H122").Value = 2132.25
$ws.Range("I122").Value = 2003.5
$ws.Range("J122").Value = 2261
$ws.Range("K122").Value = 6010.5
$ws.Range("L122").Value = 6783
$ws.Range("M122").Value = -3560.5
$ws.Range("N122").Value = -11683
$ws.Range("H132").Value = 181827710
$ws.Range("I132").Value = 666667140
$ws.Range("J132").Value = 12928.125
$ws.Range("K132").Value = 2000001420
$ws.Range("L132").Value = 38784.375
$ws.Range("M132").Value = -1999998890
$ws.Range("N132").Value = -43844.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6411840
$ws.Range("I7").Value = 1976
$ws.Range("K7").Value = 1976
$ws.Range("M7").Value = -1864
$ws.Range("H98").Value = 23700
$ws.Range("J98").Value = 23700
$ws.Range("L98").Value = 23700
$ws.Range("N98").Value = -29690
$ws.Range("H126").Value = 6411840
$ws.Range("I126").Value = 1976
$ws.Range("K126").Value = 5928
$ws.Range("M126").Value = -3458
$ws.Range("H132").Value = 20457412
$ws.Range("I132").Value = 66668164
$ws.Range("J132").Value = 5053827
$ws.Range("K132").Value = 200004492
$ws.Range("L132").Value = 15161481
$ws.Range("M132").Value = -200001962
$ws.Range("N132").Value = -15166541
$ws.Range("H136").Value = 1087674
$ws.Range("I136").Value = 7269.1904
$ws.Range("J136").Value = 2978382.5
$ws.Range("K136").Value = 21807.5712
$ws.Range("L136").Value = 8935147.5
$ws.Range("M136").Value = -19257.5712
$ws.Range("N136").Value = -8940247.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 155826880
$ws.Range("I132").Value = 240007800
$ws.Range("J132").Value = 85676110
$ws.Range("K132").Value = 720023400
$ws.Range("L132").Value = 257028330
$ws.Range("M132").Value = -720020870
$ws.Range("N132").Value = -257033390
$ws.Range("H136").Value = 34288556
$ws.Range("I136").Value = 21011468
$ws.Range("J136").Value = 62502376
$ws.Range("K136").Value = 63034404
$ws.Range("L136").Value = 187507128
$ws.Range("M136").Value = -63031854
$ws.Range("N136").Value = -187512228

Write-Host "Applied all Fenrir Profits cell updates"
